$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet "AddLine" (sheet2): add the new header columns (E:K) and data rows
# ---------------------------------------------------------------------------

# -- Header row (row 1), columns E1:K1 -------------------------------------
$ws2.Range("E1").Value = "Inventory Division"
$ws2.Range("F1").Value = "Sales Division"
$ws2.Range("G1").Value = "Customer"
$ws2.Range("H1").Value = "Manufacturing User"
$ws2.Range("I1").Value = "Product"
$ws2.Range("J1").Value = "Line Number"
$ws2.Range("K1").Value = "Record Processed"

# Build the bold/grey header font once (Arial 8 FF4A4A56) by starting from the
# existing bold-Arial-8 style already used on A1:C1 and only changing the
# color - this reuses the existing font id instead of creating a brand new
# one from scratch.
$ws2.Range("A1").Copy()
$ws2.Range("E1").PasteSpecial(-4122)
$ws2.Range("E1").Font.Color = 5655114

# Re-use the freshly built style (font+xf) for the rest of the header cells.
$ws2.Range("E1").Copy()
$ws2.Range("F1:K1").PasteSpecial(-4122)

# Re-apply the values (PasteSpecial of formats only shouldn't disturb them,
# but make sure they are correct regardless).
$ws2.Range("E1").Value = "Inventory Division"
$ws2.Range("F1").Value = "Sales Division"
$ws2.Range("G1").Value = "Customer"
$ws2.Range("H1").Value = "Manufacturing User"
$ws2.Range("I1").Value = "Product"
$ws2.Range("J1").Value = "Line Number"
$ws2.Range("K1").Value = "Record Processed"

# -- Row 2 --------------------------------------------------------------
$ws2.Range("E2").Value = "a7O410000004Nwj"
$ws2.Range("F2").Value = "a7O410000004Nwj"
$ws2.Range("G2").Value = "a5B41000000PRNX"
$ws2.Range("H2").Value = "a811K000000k9cL"
$ws2.Range("I2").Value = "a6J1K000000Qgsa"
$ws2.Range("J2").Value = 1
$ws2.Range("K2").Value = $true

# -- Row 3 --------------------------------------------------------------
$ws2.Range("E3").Value = "a7O410000004Nwj"
$ws2.Range("F3").Value = "a7O410000004Nwj"
$ws2.Range("G3").Value = "a5B41000000PRNX"
$ws2.Range("H3").Value = "a811K000000k9cL"
$ws2.Range("I3").Value = "a6J1K000000Qgsf"
$ws2.Range("J3").Value = 2
$ws2.Range("K3").Value = $true

# -- Row 4 --------------------------------------------------------------
$ws2.Range("E4").Value = "a7O410000004Nwj"
$ws2.Range("F4").Value = "a7O410000004Nwj"
$ws2.Range("G4").Value = "a5B41000000PRNX"
$ws2.Range("H4").Value = "a811K000000k9cL"
$ws2.Range("I4").Value = "a6J1K000000Qgsk"
$ws2.Range("J4").Value = 3
$ws2.Range("K4").Value = $true

# -- Row 5 --------------------------------------------------------------
$ws2.Range("E5").Value = "a7O410000004Nwj"
$ws2.Range("F5").Value = "a7O410000004Nwj"
$ws2.Range("G5").Value = "a5B41000000PRNX"
$ws2.Range("H5").Value = "a811K000000k9cL"
$ws2.Range("I5").Value = "a6J1K000000Qgsp"
$ws2.Range("J5").Value = 4
$ws2.Range("K5").Value = $true

# Build the small "record" font once (Arial 9 FF222222) on I4, re-using the
# default-unstyled cell as the starting point.
$ws2.Range("I4").Font.Name = "Arial"
$ws2.Range("I4").Font.Size = 9
$ws2.Range("I4").Font.Color = 2236962

# Re-use that style for I5.
$ws2.Range("I4").Copy()
$ws2.Range("I5").PasteSpecial(-4122)
$ws2.Range("I5").Value = "a6J1K000000Qgsp"

# -- Column width for the new "Product" (I) column ------------------------
$ws2.Range("I1").Value = "Product"
$ws2.Columns.Item(9).ColumnWidth = 15.6

# -- Selection -------------------------------------------------------------
$ws2.Activate()
$ws2.Range("D14").Select()

# ---------------------------------------------------------------------------
# Sheet "AddHeader" (sheet1): only the cursor/selection moved
# ---------------------------------------------------------------------------
$ws1.Range("F15").Select()

# Restore AddLine as the active / visible sheet (it was already the active
# tab before these edits, keep it that way).
$ws2.Activate()
$ws2.Range("D14").Select()
